# FuelPrices update (2025-04-25 02:47):
#  - Column A now holds the Date values (previously held in column C) and
#    column C now holds the MLBSO00 values (previously held in column A).
#    i.e. the A and C columns are swapped, including the header labels.
#  - A new row of data (today's reading) is inserted right after the header:
#       Date = 45768, LNBSF00 = 758.1079999999999, MLBSO00 = 755.163

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1   # last row with original data (27)

# --- capture the original data (rows 2..lastRow) before mutating anything ---
$oldDate = @{}   # previously column C
$oldLN   = @{}   # column B (LNBSF00, unaffected by the swap)
$oldML   = @{}   # previously column A (MLBSO00)
for ($r = 2; $r -le $lastRow; $r++) {
    $oldDate[$r] = $ws.Cells.Item($r, 3).Value2
    $oldLN[$r]   = $ws.Cells.Item($r, 2).Value2
    $oldML[$r]   = $ws.Cells.Item($r, 1).Value2
}

# --- header row: swap A1/C1 labels, B1 stays as-is ---
$ws.Cells.Item(1, 1).Value = "Date"
$ws.Cells.Item(1, 3).Value = "MLBSO00"

# --- rewrite rows 3..lastRow+1 from the captured old rows 2..lastRow, ---
# --- effectively inserting a new row directly under the header.        ---
# Work from the bottom up so source rows are read before being overwritten.
for ($r = $lastRow; $r -ge 2; $r--) {
    $dest = $r + 1

    $ws.Cells.Item($dest, 1).Style = "Normal"
    $ws.Cells.Item($dest, 1).Value = $oldDate[$r]
    $ws.Cells.Item($dest, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($dest, 2).Style = "Normal"
    if ($oldLN[$r] -eq "") {
        $ws.Cells.Item($dest, 2).Value = $null
    } else {
        $ws.Cells.Item($dest, 2).Value = $oldLN[$r]
    }

    $ws.Cells.Item($dest, 3).Style = "Normal"
    if ($oldML[$r] -eq "") {
        $ws.Cells.Item($dest, 3).Value = $null
    } else {
        $ws.Cells.Item($dest, 3).Value = $oldML[$r]
    }
}

# --- new row 2: latest reading ---
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 1).Value = 45768
$ws.Cells.Item(2, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(2, 2).Style = "Normal"
$ws.Cells.Item(2, 2).Value = 758.1079999999999

$ws.Cells.Item(2, 3).Style = "Normal"
$ws.Cells.Item(2, 3).Value = 755.163

Write-Host "FuelPrices sheet updated: inserted latest row and swapped Date/MLBSO00 columns."
